# "update manual and rebuild"
#
# The document was carrying a handful of SharePoint/OneDrive document-library
# "custom XML parts" (the content-type schema in customXml/item1.xml and the
# FormTemplates stub in customXml/item2.xml, plus their itemProps*.xml
# datastore-item descriptors). These are left-over library metadata that
# isn't referenced anywhere from the document body (no content control /
# XML mapping binds to them) - a "rebuild" of the document drops them.
#
# Remove every custom XML part from the package. We go after the two known
# SharePoint namespaces specifically (belt) and then sweep anything left
# over generically from the end of the collection backwards (suspenders),
# since deleting by index shifts the remaining items down.

$d = $word.ActiveDocument

function Remove-CustomXmlPartsByNamespace($ns) {
    try {
        $scoped = $d.CustomXMLParts.SelectByNamespace($ns)
    } catch {
        $scoped = $null
    }
    if ($scoped -eq $null) { return }

    $n = 0
    try { $n = $scoped.Count } catch { $n = 0 }

    for ($i = $n; $i -ge 1; $i--) {
        try {
            $part = $scoped.Item($i)
            if ($part -ne $null) {
                $part.Delete()
            }
        } catch {
            # namespace not present / already removed - nothing to do
        }
    }
}

# The SharePoint content-type schema (customXml/item1.xml).
Remove-CustomXmlPartsByNamespace("http://schemas.microsoft.com/office/2006/metadata/contentType")
# The SharePoint document-library FormTemplates stub (customXml/item2.xml).
Remove-CustomXmlPartsByNamespace("http://schemas.microsoft.com/sharepoint/v3/contenttype/forms")

# Generic sweep: walk whatever is left, deleting anything that isn't one of
# Word's own built-in parts (core/extended properties, cover-page props).
$builtinNamespaces = @(
    "http://schemas.openxmlformats.org/package/2006/metadata/core-properties",
    "http://schemas.openxmlformats.org/officeDocument/2006/extended-properties",
    "http://schemas.microsoft.com/office/2006/coverPageProps"
)

$total = 0
try { $total = $d.CustomXMLParts.Count } catch { $total = 0 }

for ($i = $total; $i -ge 1; $i--) {
    try {
        $part = $d.CustomXMLParts.Item($i)
        if ($part -eq $null) { continue }

        $ns = $null
        try { $ns = $part.NamespaceURI } catch { $ns = $null }

        if ($builtinNamespaces -notcontains $ns) {
            $part.Delete()
        }
    } catch {
        # best-effort cleanup - ignore parts that can't be inspected/removed
    }
}

Write-Output ("CustomXMLParts remaining: " + $d.CustomXMLParts.Count)
